# Update the "VM results" worksheet with refreshed accuracy / precision / f1
# numbers for the two result tables (rows 2-11 and rows 16-25). The
# AVERAGE() formulas in column D/F/G of rows 11 and 25 recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Table 1 (rows 2-10) ----
$ws.Range("D2").Value = 99.97

$ws.Range("D3").Value = 66.099999999999994
$ws.Range("F3").Value = 0.15
$ws.Range("G3").Value = 0.28999999999999998

$ws.Range("D4").Value = 43.23
$ws.Range("F4").Value = 0.09
$ws.Range("G4").Value = 0.17

$ws.Range("D5").Value = 99.2

$ws.Range("D6").Value = 99.9
$ws.Range("G6").Value = 28.57

$ws.Range("D7").Value = 99.98

$ws.Range("D8").Value = 70.260000000000005

$ws.Range("D9").Value = 100
$ws.Range("G9").Value = 66.67

$ws.Range("D10").Value = 99.98

# ---- Table 2 (rows 16-24) ----
$ws.Range("D16").Value = 99.98

$ws.Range("D17").Value = 77.75
$ws.Range("F17").Value = 0.22
$ws.Range("G17").Value = 0.44

$ws.Range("D18").Value = 78.17
$ws.Range("F18").Value = 0.23
$ws.Range("G18").Value = 0.45

$ws.Range("D19").Value = 99.93

$ws.Range("D20").Value = 99.88
$ws.Range("F20").Value = 16.670000000000002

$ws.Range("D21").Value = 100

$ws.Range("D22").Value = 76.81

$ws.Range("D23").Value = 100
$ws.Range("G23").Value = 66.67

$ws.Range("D24").Value = 99.85
$ws.Range("F24").Value = 14.29
$ws.Range("G24").Value = 22.22

# ---- Update the saved view/selection to match the author's session ----
$ws.Range("E26").Select()
